$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the summary/header figures ---
# Total "VALOR MORA" (overdue amount) grows because a new period was added
$ws.Range("E11").Value = 569400

# "Cant. Periodos" (number of periods) increases from 4 to 5
$ws.Range("F13").Value = 5

# --- Center-align the "Periodo Mora" column for all existing data rows ---
$ws.Range("E16:E23").HorizontalAlignment = -4108

# --- Insert two new rows just before the current last data row (23), ---
# --- duplicating its formatting so the new rows blend into the table, ---
# --- and leaving the (previously last) row 23 with the regular (non-bottom) style ---
$ws.Rows("23:24").Insert(-4121)

$ws.Range("B22:J22").Copy()
$ws.Range("B23:J24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Fill in the new rows with the new period (2509) records ---
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1047374292"
$ws.Range("D23").Value = "AMAURI BARON PITALUA"
$ws.Range("E23").Value = "2509"
$ws.Range("F23").Value = 56940
$ws.Range("G23").Value = 1423500

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1049941792"
$ws.Range("D24").Value = "CARLOS JOSE ARIAS RAMOS"
$ws.Range("E24").Value = "2509"
$ws.Range("F24").Value = 56940
$ws.Range("G24").Value = 1423500

# Keep the "Periodo Mora" column centered on the two new rows as well
$ws.Range("E23:E24").HorizontalAlignment = -4108
